$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "organization"
$ws.Range("A2").Value = '{"dunsControlStatus":{"detailedOperatingStatus":{},"subjectHandlingDetails":[],"isSelfRequestedDUNS":null,"isDelisted":false,"isMailUndeliverable":false,"firstReportDate":null,"operatingSubStatus":{"description":"Active","dnbCode":9074,"startDate":null},"isMarketable":true,"isTelephoneDisconnected":null,"operatingStatus":{"description":"Active","dnbCode":9074,"startDate":null},"selfRequestDate":null,"recordClass":{}},"websiteAddress":[],"industrialPlantsCount":null,"stockExchanges":[],"globalUltimate":{},"registeredAddress":{},"multilingualPrimaryAddress":[],"isAgent":null,"multilingualTradestyleNames":[],"defaultCurrency":"INR","employerDesignation":{},"incorporatedDate":"2021-05-25","domesticUltimate":{},"isSmallBusiness":null,"certifiedEmail":null,"multilingualRegisteredAddress":[],"financials":[{"informationScopeDnBCode":null,"informationScopeDescription":null,"financialStatementToDate":null,"yearlyRevenue":[{"currency":"INR","value":198168709},{"currency":"USD","value":2360982}],"reliabilityDnBCode":9094,"unitCode":"SingleUnits","financialStatementDuration":null,"reliabilityDescription":"Modelled"}],"registeredDetails":{"legalForm":{"registrationStatus":{},"description":"Private Limited Liability Company","dnbCode":12903,"localLegalForms":[]}},"businessTrustIndex":{},"multiLingualSearchNames":[],"standardizedStockExchanges":[],"banks":[],"telephone":[],"registrationNumbers":[{"registrationLocation":null,"registrationNumber":"AAGCE4803L","typeDnBCode":12897,"registrationNumberClass":{"description":"Fiscal / Tax Registration Number","dnbCode":2929},"isPreferredRegistrationNumber":null,"typeDescription":"Tax Identification Number (IN)"},{"registrationLocation":null,"registrationNumber":"U72900KA2021FTC147805","typeDnBCode":1362,"registrationNumberClass":{"description":"National Business Registration Number","dnbCode":2888},"isPreferredRegistrationNumber":null,"typeDescription":"Company Incorporation Number (IN)"}],"numberOfEmployees":[{"employeeFiguresDate":"2023-07-31","informationScopeDnBCode":9067,"informationScopeDescription":"Consolidated","trend":[],"reliabilityDnBCode":9092,"employeeCategories":[],"value":26,"reliabilityDescription":"Actual"},{"employeeFiguresDate":"2023-07-31","informationScopeDnBCode":9066,"informationScopeDescription":"Individual","reliabilityDnBCode":9092,"employeeCategories":[],"value":26,"reliabilityDescription":"Actual"}],"tradeStyleNames":[],"registeredName":"EXOSTAR INDIA PRIVATE LIMITED","duns":"854404636","primaryAddress":{"addressCountry":{"name":"India","isoAlpha2Code":"IN"},"postalCodePosition":{"description":"Post Code presented after Town/City name","dnbCode":1011},"geographicalPrecision":{"description":"Postal Code Centroid","dnbCode":30258},"minorTownName":null,"streetNumber":null,"continentalRegion":{"name":"Asia Pacific"},"postalCode":"560010","isRegisteredAddress":true,"latitude":12.99877,"language":{"description":"English","dnbCode":39},"addressCounty":{},"postalRoute":null,"streetName":null,"postOfficeBox":{},"isManufacturingLocation":null,"standardAddressCodes":[],"streetAddress":{"line2":"Dr. Rajkumar Road, Rajajinagar","line1":"Regus, Tejas Arcade #527/b, 3rd & 4th Floor 1st Main Road"},"statisticalArea":{},"addressLocality":{"name":"Bengaluru"},"addressRegion":{"isoSubDivisionCode":"IN-KA","name":"Karnataka","administrativeDivisionCode":null,"isoSubDivisionName":"Karnataka","abbreviatedName":"KA"},"longitude":77.57331},"startDate":"2021-05-25","primaryIndustryCode":{"usSicV4":"7379","usSicV4Description":"Computer related services"},"isExporter":null,"preferredLanguage":{"description":"English","dnbCode":39},"isFortune1000Listed":null,"controlOwnershipDate":"2021-05-25","industryCodes":[{"code":"541512","typeDnBCode":37788,"description":"Computer Systems Design Services","typeDescription":"North American Industry Classification System 2022","priority":1},{"code":"73790000","typeDnBCode":3599,"description":"Computer related services, nec","typeDescription":"D&B Standard Industry Code","priority":1},{"code":"6203","typeDnBCode":29104,"description":"Computer facilities management activities","typeDescription":"NACE Revision 2","priority":1},{"code":"7379","typeDnBCode":399,"description":"Computer related services","typeDescription":"US Standard Industry Code 1987 - 4 digit","priority":1},{"code":"42","typeDnBCode":35912,"description":"Computer System Design Services","typeDescription":"D&B Hoovers Industry Classification","priority":1},{"code":"I","typeDnBCode":24657,"description":"Services","typeDescription":"D&B Standard Major Industry Code","priority":1},{"code":"6202","typeDnBCode":42726,"description":"Computer consultancy and computer facilities management activities","typeDescription":"ISIC Revision 4","priority":1}],"imperialCalendarStartYear":null,"countryISOAlpha2Code":"IN","multilingualPrimaryName":[],"fiscalYearEnd":null,"organizationSizeCategory":{},"corporateLinkage":{},"isImporter":null,"unspscCodes":[{"code":"81110000","description":"Computer services","priority":1}],"isStandalone":true,"isForbesLargestPrivateCompaniesListed":null,"tsrReportDate":null,"controlOwnershipType":{},"email":[{"address":"sudeep.dharan@exostar.com"}],"summary":[],"legalEntityIdentifier":null,"securitiesReportID":null,"businessEntityType":{"description":"Corporation","dnbCode":451},"primaryName":"EXOSTAR INDIA PRIVATE LIMITED","legalForm":{"registrationLocation":{},"description":"Private Limited Company","dnbCode":31148,"startDate":"2021-05-25"},"charterType":{},"mailingAddress":{},"activities":[],"isNonClassifiedEstablishment":null,"multilingualRegisteredNames":[],"subjectComments":[],"investigationDate":null}'

$ws.Range("B1:F2").EntireColumn.Delete() | Out-Null

# Excel auto-sizes column A to fit the long JSON value (bestFit), matching the
# target workbook's column width of ~255.63 characters for this content.
$ws.Columns.Item(1).ColumnWidth = 254.8

$ws.Range("A2").Select() | Out-Null
